$wb = $excel.ActiveWorkbook

# ----- Sheet ALC: updating 2 row(s) of recalculated price/profit data -----
$ws = $wb.Worksheets.Item("ALC")
# Row 138
$ws.Cells.Item(138, 8).Value = 0  # H138: 2000 -> 0
$ws.Cells.Item(138, 9).Value = 0  # I138: 2000 -> 0
$ws.Cells.Item(138, 11).Value = 0  # K138: 6000 -> 0
$ws.Cells.Item(138, 13).ClearContents()  # M138: -860 -> (removed)
# Row 141
$ws.Cells.Item(141, 8).Value = 2075.1667  # H141: 2227.5454 -> 2075.1667
$ws.Cells.Item(141, 9).Value = 1145  # I141: 1238.25 -> 1145
$ws.Cells.Item(141, 11).Value = 3435  # K141: 3714.75 -> 3435
$ws.Cells.Item(141, 13).Value = 1745  # M141: 1465.25 -> 1745

# ----- Sheet ARM: updating 7 row(s) of recalculated price/profit data -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 2228.5652  # H2: 2539.65 -> 2228.5652
$ws.Cells.Item(2, 9).Value = 1625.9445  # I2: 1920.2 -> 1625.9445
$ws.Cells.Item(2, 11).Value = 1625.9445  # K2: 1920.2 -> 1625.9445
$ws.Cells.Item(2, 13).Value = -1512.9445  # M2: -1807.2 -> -1512.9445
# Row 32
$ws.Cells.Item(32, 8).Value = 4311.9473  # H32: 4343.5264 -> 4311.9473
$ws.Cells.Item(32, 9).Value = 4311.9473  # I32: 4343.5264 -> 4311.9473
$ws.Cells.Item(32, 11).Value = 4311.9473  # K32: 4343.5264 -> 4311.9473
$ws.Cells.Item(32, 13).Value = -4024.9473  # M32: -4056.5264 -> -4024.9473
# Row 74
$ws.Cells.Item(74, 8).Value = 1242.2142  # H74: 1416.7273 -> 1242.2142
$ws.Cells.Item(74, 9).Value = 1242.2142  # I74: 1416.7273 -> 1242.2142
$ws.Cells.Item(74, 11).Value = 1242.2142  # K74: 1416.7273 -> 1242.2142
$ws.Cells.Item(74, 13).Value = -368.2141999999999  # M74: -542.7273 -> -368.2141999999999
# Row 77
$ws.Cells.Item(77, 8).Value = 1242.2142  # H77: 1416.7273 -> 1242.2142
$ws.Cells.Item(77, 9).Value = 1242.2142  # I77: 1416.7273 -> 1242.2142
$ws.Cells.Item(77, 11).Value = 6211.071  # K77: 7083.636500000001 -> 6211.071
$ws.Cells.Item(77, 13).Value = -1843.071  # M77: -2715.636500000001 -> -1843.071
# Row 93
$ws.Cells.Item(93, 8).Value = 49998  # H93: 11111 -> 49998
$ws.Cells.Item(93, 9).Value = 49998  # I93: 0 -> 49998
$ws.Cells.Item(93, 10).Value = 0  # J93: 11111 -> 0
$ws.Cells.Item(93, 11).Value = 49998  # K93: 0 -> 49998
$ws.Cells.Item(93, 12).Value = 0  # L93: 11111 -> 0
$ws.Cells.Item(93, 13).Value = -47502  # M93: None -> -47502
$ws.Cells.Item(93, 14).ClearContents()  # N93: -16103 -> (removed)
# Row 116
$ws.Cells.Item(116, 8).Value = 2228.5652  # H116: 2539.65 -> 2228.5652
$ws.Cells.Item(116, 9).Value = 1625.9445  # I116: 1920.2 -> 1625.9445
$ws.Cells.Item(116, 11).Value = 1625.9445  # K116: 1920.2 -> 1625.9445
$ws.Cells.Item(116, 13).Value = 668.0554999999999  # M116: 373.8 -> 668.0554999999999
# Row 122
$ws.Cells.Item(122, 8).Value = 1648.8  # H122: 1729 -> 1648.8
$ws.Cells.Item(122, 9).Value = 1648.8  # I122: 1729 -> 1648.8
$ws.Cells.Item(122, 11).Value = 4946.4  # K122: 5187 -> 4946.4
$ws.Cells.Item(122, 13).Value = -2496.4  # M122: -2737 -> -2496.4

# ----- Sheet BSM: updating 4 row(s) of recalculated price/profit data -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 2228.5652  # H3: 2539.65 -> 2228.5652
$ws.Cells.Item(3, 9).Value = 1625.9445  # I3: 1920.2 -> 1625.9445
$ws.Cells.Item(3, 11).Value = 1625.9445  # K3: 1920.2 -> 1625.9445
$ws.Cells.Item(3, 13).Value = -1511.9445  # M3: -1806.2 -> -1511.9445
# Row 13
$ws.Cells.Item(13, 8).Value = 59900  # H13: 58949.5 -> 59900
$ws.Cells.Item(13, 10).Value = 59900  # J13: 58949.5 -> 59900
$ws.Cells.Item(13, 12).Value = 59900  # L13: 58949.5 -> 59900
$ws.Cells.Item(13, 14).Value = -60236  # N13: -59285.5 -> -60236
# Row 105
$ws.Cells.Item(105, 8).Value = 2080.8  # H105: 2073 -> 2080.8
$ws.Cells.Item(105, 9).Value = 2074.625  # I105: 2065.7778 -> 2074.625
$ws.Cells.Item(105, 11).Value = 2074.625  # K105: 2065.7778 -> 2074.625
$ws.Cells.Item(105, 13).Value = -327.625  # M105: -318.7777999999998 -> -327.625
# Row 134
$ws.Cells.Item(134, 8).Value = 2892.2778  # H134: 3429.394 -> 2892.2778
$ws.Cells.Item(134, 9).Value = 2689.2  # I134: 3224.0625 -> 2689.2
$ws.Cells.Item(134, 11).Value = 8067.599999999999  # K134: 9672.1875 -> 8067.599999999999
$ws.Cells.Item(134, 13).Value = -5532.599999999999  # M134: -7137.1875 -> -5532.599999999999

# ----- Sheet CRP: updating 2 row(s) of recalculated price/profit data -----
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Cells.Item(6, 8).Value = 8677.909  # H6: 1002 -> 8677.909
$ws.Cells.Item(6, 9).Value = 383.8889  # I6: 0 -> 383.8889
$ws.Cells.Item(6, 10).Value = 46001  # J6: 1002 -> 46001
$ws.Cells.Item(6, 11).Value = 383.8889  # K6: 0 -> 383.8889
$ws.Cells.Item(6, 12).Value = 46001  # L6: 1002 -> 46001
$ws.Cells.Item(6, 13).Value = -270.8889  # M6: None -> -270.8889
$ws.Cells.Item(6, 14).Value = -46227  # N6: -1228 -> -46227
# Row 122
$ws.Cells.Item(122, 8).Value = 1745.7222  # H122: 1772.1765 -> 1745.7222
$ws.Cells.Item(122, 9).Value = 1672.8235  # I122: 1696.375 -> 1672.8235
$ws.Cells.Item(122, 11).Value = 5018.470499999999  # K122: 5089.125 -> 5018.470499999999
$ws.Cells.Item(122, 13).Value = -2568.470499999999  # M122: -2639.125 -> -2568.470499999999

# ----- Sheet CUL: updating 3 row(s) of recalculated price/profit data -----
$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Cells.Item(107, 8).Value = 297.25  # H107: 409.8 -> 297.25
$ws.Cells.Item(107, 9).Value = 200  # I107: 525 -> 200
$ws.Cells.Item(107, 10).Value = 329.66666  # J107: 333 -> 329.66666
$ws.Cells.Item(107, 11).Value = 600  # K107: 1575 -> 600
$ws.Cells.Item(107, 12).Value = 988.9999799999999  # L107: 999 -> 988.9999799999999
$ws.Cells.Item(107, 13).Value = 1320  # M107: 345 -> 1320
$ws.Cells.Item(107, 14).Value = -4828.99998  # N107: -4839 -> -4828.99998
# Row 131
$ws.Cells.Item(131, 8).Value = 997.5  # H131: 1000 -> 997.5
$ws.Cells.Item(131, 9).Value = 997.5  # I131: 1000 -> 997.5
$ws.Cells.Item(131, 10).Value = 0  # J131: 1000 -> 0
$ws.Cells.Item(131, 11).Value = 2992.5  # K131: 3000 -> 2992.5
$ws.Cells.Item(131, 12).Value = 0  # L131: 3000 -> 0
$ws.Cells.Item(131, 13).Value = 2047.5  # M131: 2040 -> 2047.5
$ws.Cells.Item(131, 14).ClearContents()  # N131: -13080 -> (removed)
# Row 140
$ws.Cells.Item(140, 8).Value = 1723.3  # H140: 1731.85 -> 1723.3
$ws.Cells.Item(140, 9).Value = 1150.9412  # I140: 1161 -> 1150.9412
$ws.Cells.Item(140, 11).Value = 3452.8236  # K140: 3483 -> 3452.8236
$ws.Cells.Item(140, 13).Value = 1727.1764  # M140: 1697 -> 1727.1764

# ----- Sheet GSM: updating 5 row(s) of recalculated price/profit data -----
$ws = $wb.Worksheets.Item("GSM")
# Row 26
$ws.Cells.Item(26, 8).Value = 0  # H26: 45000 -> 0
$ws.Cells.Item(26, 10).Value = 0  # J26: 45000 -> 0
$ws.Cells.Item(26, 12).Value = 0  # L26: 45000 -> 0
$ws.Cells.Item(26, 14).ClearContents()  # N26: -45560 -> (removed)
# Row 50
$ws.Cells.Item(50, 8).Value = 0  # H50: 45000 -> 0
$ws.Cells.Item(50, 10).Value = 0  # J50: 45000 -> 0
$ws.Cells.Item(50, 12).Value = 0  # L50: 45000 -> 0
$ws.Cells.Item(50, 14).ClearContents()  # N50: -45996 -> (removed)
# Row 102
$ws.Cells.Item(102, 8).Value = 1860.3793  # H102: 1925.0714 -> 1860.3793
$ws.Cells.Item(102, 9).Value = 1498.1538  # I102: 1556.12 -> 1498.1538
$ws.Cells.Item(102, 11).Value = 1498.1538  # K102: 1556.12 -> 1498.1538
$ws.Cells.Item(102, 13).Value = 123.8462  # M102: 65.88000000000011 -> 123.8462
# Row 113
$ws.Cells.Item(113, 8).Value = 9949  # H113: 9999.166999999999 -> 9949
$ws.Cells.Item(113, 9).Value = 9750  # I113: 0 -> 9750
$ws.Cells.Item(113, 10).Value = 9998.75  # J113: 9999.166999999999 -> 9998.75
$ws.Cells.Item(113, 11).Value = 9750  # K113: 0 -> 9750
$ws.Cells.Item(113, 12).Value = 9998.75  # L113: 9999.166999999999 -> 9998.75
$ws.Cells.Item(113, 13).Value = -7580  # M113: None -> -7580
$ws.Cells.Item(113, 14).Value = -14338.75  # N113: -14339.167 -> -14338.75
# Row 126
$ws.Cells.Item(126, 8).Value = 3737.2  # H126: 3810.5334 -> 3737.2
$ws.Cells.Item(126, 9).Value = 3466  # I126: 3550.6155 -> 3466
$ws.Cells.Item(126, 11).Value = 10398  # K126: 10651.8465 -> 10398
$ws.Cells.Item(126, 13).Value = -7928  # M126: -8181.8465 -> -7928

# ----- Sheet LTW: updating 10 row(s) of recalculated price/profit data -----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 2189.1667  # H22: 2858 -> 2189.1667
$ws.Cells.Item(22, 9).Value = 877.6667  # I22: 930 -> 877.6667
$ws.Cells.Item(22, 11).Value = 877.6667  # K22: 930 -> 877.6667
$ws.Cells.Item(22, 13).Value = -582.6667  # M22: -635 -> -582.6667
# Row 27
$ws.Cells.Item(27, 8).Value = 2189.1667  # H27: 2858 -> 2189.1667
$ws.Cells.Item(27, 9).Value = 877.6667  # I27: 930 -> 877.6667
$ws.Cells.Item(27, 11).Value = 877.6667  # K27: 930 -> 877.6667
$ws.Cells.Item(27, 13).Value = -770.6667  # M27: -823 -> -770.6667
# Row 30
$ws.Cells.Item(30, 8).Value = 285.25  # H30: 450 -> 285.25
$ws.Cells.Item(30, 9).Value = 285.25  # I30: 450 -> 285.25
$ws.Cells.Item(30, 11).Value = 285.25  # K30: 450 -> 285.25
$ws.Cells.Item(30, 13).Value = -177.25  # M30: -342 -> -177.25
# Row 55
$ws.Cells.Item(55, 8).Value = 1101.7916  # H55: 1134.0435 -> 1101.7916
$ws.Cells.Item(55, 9).Value = 1582.1666  # I55: 1826.6 -> 1582.1666
$ws.Cells.Item(55, 11).Value = 1582.1666  # K55: 1826.6 -> 1582.1666
$ws.Cells.Item(55, 13).Value = -1409.1666  # M55: -1653.6 -> -1409.1666
# Row 82
$ws.Cells.Item(82, 8).Value = 5138.2856  # H82: 5264.385 -> 5138.2856
$ws.Cells.Item(82, 10).Value = 5855.875  # J82: 6192.5713 -> 5855.875
$ws.Cells.Item(82, 12).Value = 5855.875  # L82: 6192.5713 -> 5855.875
$ws.Cells.Item(82, 14).Value = -6577.875  # N82: -6914.5713 -> -6577.875
# Row 85
$ws.Cells.Item(85, 8).Value = 5138.2856  # H85: 5264.385 -> 5138.2856
$ws.Cells.Item(85, 10).Value = 5855.875  # J85: 6192.5713 -> 5855.875
$ws.Cells.Item(85, 12).Value = 5855.875  # L85: 6192.5713 -> 5855.875
$ws.Cells.Item(85, 14).Value = -8351.875  # N85: -8688.5713 -> -8351.875
# Row 100
$ws.Cells.Item(100, 8).Value = 5097.6523  # H100: 5268.0454 -> 5097.6523
$ws.Cells.Item(100, 9).Value = 1976.909  # I100: 2039.7 -> 1976.909
$ws.Cells.Item(100, 11).Value = 1976.909  # K100: 2039.7 -> 1976.909
$ws.Cells.Item(100, 13).Value = -1435.909  # M100: -1498.7 -> -1435.909
# Row 122
$ws.Cells.Item(122, 8).Value = 3989.25  # H122: 3991 -> 3989.25
$ws.Cells.Item(122, 9).Value = 3989.25  # I122: 3991 -> 3989.25
$ws.Cells.Item(122, 11).Value = 11967.75  # K122: 11973 -> 11967.75
$ws.Cells.Item(122, 13).Value = -9517.75  # M122: -9523 -> -9517.75
# Row 132
$ws.Cells.Item(132, 8).Value = 0  # H132: 3444.5 -> 0
$ws.Cells.Item(132, 9).Value = 0  # I132: 3444.5 -> 0
$ws.Cells.Item(132, 11).Value = 0  # K132: 10333.5 -> 0
$ws.Cells.Item(132, 13).ClearContents()  # M132: -7803.5 -> (removed)
# Row 136
$ws.Cells.Item(136, 8).Value = 5331.3335  # H136: 5499.75 -> 5331.3335
$ws.Cells.Item(136, 9).Value = 4997.25  # I136: 5000 -> 4997.25
$ws.Cells.Item(136, 11).Value = 14991.75  # K136: 15000 -> 14991.75
$ws.Cells.Item(136, 13).Value = -12441.75  # M136: -12450 -> -12441.75

# ----- Sheet WVR: updating 11 row(s) of recalculated price/profit data -----
$ws = $wb.Worksheets.Item("WVR")
# Row 21
$ws.Cells.Item(21, 8).Value = 2000000  # H21: 2750000 -> 2000000
$ws.Cells.Item(21, 9).Value = 2000000  # I21: 2750000 -> 2000000
$ws.Cells.Item(21, 11).Value = 2000000  # K21: 2750000 -> 2000000
$ws.Cells.Item(21, 13).Value = -1999765  # M21: -2749765 -> -1999765
# Row 25
$ws.Cells.Item(25, 8).Value = 41513.5  # H25: 63027 -> 41513.5
$ws.Cells.Item(25, 10).Value = 41513.5  # J25: 63027 -> 41513.5
$ws.Cells.Item(25, 12).Value = 41513.5  # L25: 63027 -> 41513.5
$ws.Cells.Item(25, 14).Value = -42099.5  # N25: -63613 -> -42099.5
# Row 28
$ws.Cells.Item(28, 8).Value = 20000  # H28: 0 -> 20000
$ws.Cells.Item(28, 10).Value = 20000  # J28: 0 -> 20000
$ws.Cells.Item(28, 12).Value = 20000  # L28: 0 -> 20000
$ws.Cells.Item(28, 14).Value = -20696  # N28: None -> -20696
# Row 29
$ws.Cells.Item(29, 8).Value = 10000  # H29: 5000 -> 10000
$ws.Cells.Item(29, 9).Value = 0  # I29: 5000 -> 0
$ws.Cells.Item(29, 10).Value = 10000  # J29: 0 -> 10000
$ws.Cells.Item(29, 11).Value = 0  # K29: 5000 -> 0
$ws.Cells.Item(29, 12).Value = 10000  # L29: 0 -> 10000
$ws.Cells.Item(29, 13).ClearContents()  # M29: -4710 -> (removed)
$ws.Cells.Item(29, 14).Value = -10580  # N29: None -> -10580
# Row 30
$ws.Cells.Item(30, 8).Value = 41005  # H30: 57010 -> 41005
$ws.Cells.Item(30, 9).Value = 25000  # I30: 0 -> 25000
$ws.Cells.Item(30, 11).Value = 25000  # K30: 0 -> 25000
$ws.Cells.Item(30, 13).Value = -24893  # M30: None -> -24893
# Row 35
$ws.Cells.Item(35, 8).Value = 2000000  # H35: 2750000 -> 2000000
$ws.Cells.Item(35, 9).Value = 2000000  # I35: 2750000 -> 2000000
$ws.Cells.Item(35, 11).Value = 2000000  # K35: 2750000 -> 2000000
$ws.Cells.Item(35, 13).Value = -1999710  # M35: -2749710 -> -1999710
# Row 82
$ws.Cells.Item(82, 8).Value = 41250  # H82: 40000 -> 41250
$ws.Cells.Item(82, 10).Value = 41250  # J82: 40000 -> 41250
$ws.Cells.Item(82, 12).Value = 41250  # L82: 40000 -> 41250
$ws.Cells.Item(82, 14).Value = -42016  # N82: -40766 -> -42016
# Row 85
$ws.Cells.Item(85, 8).Value = 41250  # H85: 40000 -> 41250
$ws.Cells.Item(85, 10).Value = 41250  # J85: 40000 -> 41250
$ws.Cells.Item(85, 12).Value = 41250  # L85: 40000 -> 41250
$ws.Cells.Item(85, 14).Value = -43902  # N85: -42652 -> -43902
# Row 107
$ws.Cells.Item(107, 8).Value = 570.6429000000001  # H107: 650.2308 -> 570.6429000000001
$ws.Cells.Item(107, 9).Value = 528.25  # I107: 557.3 -> 528.25
$ws.Cells.Item(107, 10).Value = 825  # J107: 960 -> 825
$ws.Cells.Item(107, 11).Value = 1584.75  # K107: 1671.9 -> 1584.75
$ws.Cells.Item(107, 12).Value = 2475  # L107: 2880 -> 2475
$ws.Cells.Item(107, 13).Value = 335.25  # M107: 248.1000000000001 -> 335.25
$ws.Cells.Item(107, 14).Value = -6315  # N107: -6720 -> -6315
# Row 122
$ws.Cells.Item(122, 8).Value = 1930.6  # H122: 2244 -> 1930.6
$ws.Cells.Item(122, 9).Value = 1791.3572  # I122: 2116.8462 -> 1791.3572
$ws.Cells.Item(122, 10).Value = 3880  # J122: 2381.75 -> 3880
$ws.Cells.Item(122, 11).Value = 5374.071599999999  # K122: 6350.5386 -> 5374.071599999999
$ws.Cells.Item(122, 12).Value = 11640  # L122: 7145.25 -> 11640
$ws.Cells.Item(122, 13).Value = -2924.071599999999  # M122: -3900.5386 -> -2924.071599999999
$ws.Cells.Item(122, 14).Value = -16540  # N122: -12045.25 -> -16540
# Row 132
$ws.Cells.Item(132, 8).Value = 4900.7856  # H132: 4696.3335 -> 4900.7856
$ws.Cells.Item(132, 9).Value = 5261.2  # I132: 4949.636 -> 5261.2
$ws.Cells.Item(132, 11).Value = 15783.6  # K132: 14848.908 -> 15783.6
$ws.Cells.Item(132, 13).Value = -13253.6  # M132: -12318.908 -> -13253.6

